# Fix the typo/redaction issues from the hard-copy CV review.
$d = $word.ActiveDocument

# 1. "GraphyQl" -> "GraphQl" (typo in the "currently excited about" blurb).
$d.Content.Find.Execute("GraphyQl", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "GraphQl", 2) | Out-Null

# 2. Job-heading line: "Undisclosed (name omitted) " (four runs: the bold
#    "Undisclosed (", the bold+italic "name omitted", the bold ")" and a
#    trailing plain space) collapses into a single bold run "Dataffirm ",
#    revealing the employer's real name. Matching on the full original
#    phrase (spanning all four runs) and replacing it in one go gives a
#    single resulting run that keeps the leading run's bold formatting,
#    exactly like the target edit.
$d.Content.Find.Execute("Undisclosed (name omitted) ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Dataffirm ", 2) | Out-Null

# 3. Company blurb paragraph: "Undisclosed are a financial tech startup..."
#    -> "Dataffirm are a financial tech startup...".
$d.Content.Find.Execute("Undisclosed are", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Dataffirm are", 2) | Out-Null
